$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# productName_1, productName_2, productName_5 all changed from "ADIDAS ORIGINAL" to "IPHONE 13 PRO"
$ws.Range("B7").Value = "IPHONE 13 PRO"
$ws.Range("B11").Value = "IPHONE 13 PRO"
$ws.Range("B15").Value = "IPHONE 13 PRO"

# update the view: scroll so row 4 is at top, and select B15
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B15").Select()
